$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 248, shifting existing rows 248:274 down to 249:275.
$ws.Rows.Item(248).Insert()

# Populate the newly inserted row 248 with its data.
$ws.Cells.Item(248, 1).Value = 5
$ws.Cells.Item(248, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(248, 3).Value = "Maule"
$ws.Cells.Item(248, 4).Value = 44769
$ws.Cells.Item(248, 5).Value = 7
$ws.Cells.Item(248, 6).Value = 100112045
$ws.Cells.Item(248, 7).Value = "Zapallo"
$ws.Cells.Item(248, 8).Value = "Camote"
$ws.Cells.Item(248, 9).Value = "1a (guarda)"
$ws.Cells.Item(248, 10).Value = 900
$ws.Cells.Item(248, 11).Value = 750
$ws.Cells.Item(248, 12).Value = 750
$ws.Cells.Item(248, 13).Value = 750
$ws.Cells.Item(248, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(248, 15).Value = "Región del Maule"
$ws.Cells.Item(248, 16).Value = 750
$ws.Cells.Item(248, 17).Value = 1
$ws.Cells.Item(248, 18).Value = "Hortaliza"
